$wb = $excel.ActiveWorkbook

# Duplicate the Portugal sheet (last existing sheet) to serve as the
# template for the new Slovakia market sheet, placing the copy at the end.
$src = $wb.Worksheets.Item("Portugal")
$src.Cells.Select()
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Slovakia"

# Update the market-specific cells (new shared strings get created in the
# order the cells are written: B4 first, then B2).
$newSheet.Range("B4").Value = "NGC-2930/T3228 "
$newSheet.Range("B2").Value = "Slovakia Market"

# The shorter replacement text no longer needs the extra wrapped line that
# the Portugal template had for rows 3 and 4, so auto-fit those two rows
# back down to the default height (row 5 is left untouched).
$newSheet.Rows.Item(3).AutoFit()
$newSheet.Rows.Item(4).AutoFit()

# Leave the new sheet with the selection used while building the list of
# printers (rows 8-12) and make it the active tab.
$newSheet.Range("A8:A12").Select()
$newSheet.Activate()
